$d = $word.ActiveDocument

$find = "Campagne Constellation de Pégase 2022"
$replace = "Campagne 2022 Constellation de Pégase"

$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
